# Apply the FlashScore weekly-odds refresh for 2024-10-27.
# 1) Update a batch of odds/handicap values that were re-priced.
# 2) Remove the two USL Championship fixtures (rows 17-18) that were
#    dropped from this week's sheet; Excel's EntireRow.Delete shifts
#    nothing below them (they were the last rows) and automatically
#    shrinks the worksheet's used range / <dimension> to A1:BD16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Row 3 (Western United v WS Wanderers) --
$ws.Range("G3").Value2  = 2.7
$ws.Range("I3").Value2  = 2.6
$ws.Range("J3").Value2  = 3.1
$ws.Range("L3").Value2  = 3
$ws.Range("X3").Value2  = 17
$ws.Range("AB3").Value2 = 23
$ws.Range("AI3").Value2 = 15
$ws.Range("AM3").Value2 = 21
$ws.Range("AN3").Value2 = 5
$ws.Range("AX3").Value2 = 4.75

# -- Row 5 (Blaublitz v Oita Trinita) --
$ws.Range("M5").Value2 = 1.1
$ws.Range("N5").Value2 = 7

# -- Row 8 (Renofa Yamaguchi v Kofu) --
$ws.Range("M8").Value2 = 1.06
$ws.Range("N8").Value2 = 10
$ws.Range("O8").Value2 = 1.3
$ws.Range("P8").Value2 = 3.4
$ws.Range("Q8").Value2 = 2.03
$ws.Range("R8").Value2 = 1.83

# -- Row 9 (Tochigi SC v Shimizu S-Pulse) --
$ws.Range("G9").Value2  = 6
$ws.Range("M9").Value2  = 1.04
$ws.Range("N9").Value2  = 13
$ws.Range("Q9").Value2  = 1.85
$ws.Range("R9").Value2  = 2
$ws.Range("X9").Value2  = 34
$ws.Range("AB9").Value2 = 51
$ws.Range("AG9").Value2 = 351
$ws.Range("AM9").Value2 = 29
$ws.Range("AO9").Value2 = 34
$ws.Range("AP9").Value2 = 41
$ws.Range("AQ9").Value2 = 126
$ws.Range("AR9").Value2 = 151
$ws.Range("AS9").Value2 = 301
$ws.Range("AY9").Value2 = 7.5

# -- Row 10 (Yokohama FC v Okayama) --
$ws.Range("I10").Value2  = 4.1
$ws.Range("J10").Value2  = 2.6
$ws.Range("U10").Value2  = 1.91
$ws.Range("V10").Value2  = 1.8
$ws.Range("AA10").Value2 = 17
$ws.Range("AC10").Value2 = 8.5
$ws.Range("AG10").Value2 = 351
$ws.Range("AP10").Value2 = 23

# -- Row 15 (Pohang v Ulsan HD) --
$ws.Range("Q15").Value2 = 1.98
$ws.Range("R15").Value2 = 1.88

# -- Remove the two trailing USL Championship matches (rows 17 & 18) --
$ws.Range("A17:A18").EntireRow.Delete()
